$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows: add two new rows of Serial/Name pairs and update row 2's data
# Order of assignment matters for shared-string table insertion order.
$ws.Range("B2").Value = "leidiomar_corsini.STM"
$ws.Range("A2").Value = "ZTEGD1D29299"

$ws.Range("A3").Value = "ZTEGD2B331A9"
$ws.Range("B3").Value = "leaniroliveira"

$ws.Range("B4").Value = "laboratoriolidersantaluzia"
$ws.Range("A4").Value = "ZTEGD1E1FBFB"

# Layout tweaks
$ws.Columns.Item(2).ColumnWidth = 26.3
$ws.Range("A7").Select()

# Update the absolute path recorded by Excel for this workbook's folder
# (Path is normally computed by Excel from the file's location on disk and
# isn't writable via automation, but we record the intent here.)
$wb.Path = "C:\Users\lucas\Documents\ONU-TERMINATOR\"
